# Add season-record columns (Wins, Losses, Ties) to the roster table.
# Mirrors the author's fix: the previous scrape only captured team
# statistics, not the season win/loss/tie record, so three new columns
# are appended after the existing data (which ends at column AC) and
# populated with the team's season record for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled like the existing header row (bold, bordered,
# centered) by copying the format from the last existing header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Find the last used data row (falls back to 50 if detection fails).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 50 }

$wins = 76
$losses = 86
$ties = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
